# Update "想去人数" (F column) figures and a couple of "最低票价" (G column)
# values that changed from a numeric price to the literal text "不可售"
# (not for sale), as produced by the site's latest data refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 273
$ws1.Range("F4").Value  = 965
$ws1.Range("F6").Value  = 436
$ws1.Range("F7").Value  = 656
$ws1.Range("F8").Value  = 236
$ws1.Range("F10").Value = 5
$ws1.Range("F11").Value = 375
$ws1.Range("F12").Value = 176
$ws1.Range("F13").Value = 31
$ws1.Range("F14").Value = 757
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 1897
$ws1.Range("F17").Value = 411
$ws1.Range("F18").Value = 5503
$ws1.Range("F19").Value = 410
$ws1.Range("F20").Value = 506
$ws1.Range("F21").Value = 31
$ws1.Range("F22").Value = 72
$ws1.Range("F23").Value = 5
$ws1.Range("F24").Value = 175
$ws1.Range("F25").Value = 133

# ---------------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("G2").Value  = "不可售"
$ws2.Range("F5").Value  = 25
$ws2.Range("F13").Value = 108
$ws2.Range("F14").Value = 46

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 5416
$ws3.Range("F3").Value = 355
$ws3.Range("F4").Value = 338

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 5416
$ws4.Range("F4").Value  = 355
$ws4.Range("G5").Value  = "不可售"
$ws4.Range("F6").Value  = 338
$ws4.Range("F7").Value  = 273
$ws4.Range("F10").Value = 25
$ws4.Range("F14").Value = 965
$ws4.Range("F18").Value = 436
$ws4.Range("F19").Value = 656
$ws4.Range("F20").Value = 236
$ws4.Range("F23").Value = 5
$ws4.Range("F24").Value = 375
$ws4.Range("F25").Value = 176
$ws4.Range("F27").Value = 31
$ws4.Range("F29").Value = 757
$ws4.Range("F30").Value = 104
$ws4.Range("F31").Value = 108
$ws4.Range("F32").Value = 1897
$ws4.Range("F33").Value = 411
$ws4.Range("F34").Value = 5503
$ws4.Range("F35").Value = 46
$ws4.Range("F36").Value = 410
$ws4.Range("F37").Value = 506
$ws4.Range("F38").Value = 31
$ws4.Range("F39").Value = 72
$ws4.Range("F41").Value = 5
$ws4.Range("F42").Value = 175
$ws4.Range("F44").Value = 133
